$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing header cell A1 carries the bold/centered/bordered header
# style used throughout row 1. Stash a copy of that formatting in an
# out-of-the-way cell so we can re-apply it after the sheet is rebuilt.
$ws.Range("A1").Copy()
$ws.Range("Z100").PasteSpecial(-4122)   # xlPasteFormats

# Wipe the currently used range (values + formatting) so we can lay the
# sheet out fresh to match the new column layout / data.
$ws.Range("A1:E2").Clear()

# New header row: B=lang_code, C=code, D=name, E=descr, F=is_active
# (column A has no header).
$ws.Range("B1").Value = "lang_code"
$ws.Range("C1").Value = "code"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "descr"
$ws.Range("F1").Value = "is_active"

# Re-apply the stashed header style to the new header cells.
$ws.Range("Z100").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)

# New data rows: id, lang_code, code, name, descr, is_active
$data = @(
    @(0, "eng", "FNR", "Fingerprint",          "Finger prints of the applicant"),
    @(1, "eng", "IRS", "Iris",                 "Iris of the applicant"),
    @(2, "eng", "PHT", "Photo",                "Photo of the face of the applicant"),
    @(3, "fra", "FNR", "Empreintes digitales", "Empreintes digitales du demandeur"),
    @(4, "fra", "IRS", "Iris",                 "Iris du demandeur"),
    @(5, "fra", "PHT", "Photo",                "Photo du visage du demandeur")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $true
    $r++
}

# Column A (the numeric id column) also carries the bordered/bold style.
$ws.Range("Z100").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)

# Remove the temporary style-stash cell.
$ws.Range("Z100").Clear()
